# Update 2021 Target Depth Data with simulated/logged 2021 conference championship game stats.

$wb = $excel.ActiveWorkbook

# OFF sheet - Home row (row 2)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 276
$wsOff.Range("C2").Value = 204
$wsOff.Range("D2").Value = 80
$wsOff.Range("E2").Value = 41
$wsOff.Range("F2").Value = 3

# DEF sheet - Home row (row 2)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 268
$wsDef.Range("C2").Value = 178
$wsDef.Range("D2").Value = 68
$wsDef.Range("E2").Value = 27
$wsDef.Range("F2").Value = 8
